# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: header cell, matching the style of the other header cells (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# H2 / H3: plain numeric values, no special style (matches F2/G2, F3/G3)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
